$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.806.89"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "'2.322.33"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'301.96"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").Value = "'94.10"
$ws.Range("E6").Value = "  -3.72%  "
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.493"
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("D10").Value = "'33.95"
$ws.Range("E10").Value = "  -4.65%  "
$ws.Range("E11").Value = "  -2.07%  "
$ws.Range("D12").Value = "'18.72"
$ws.Range("E12").Value = "  -4.17%  "
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("D14").Value = "'6.70"
$ws.Range("E14").Value = "  -3.34%  "
$ws.Range("D15").Value = "'2.684.77"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "'2.319.18"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").Value = "'42.747.76"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").Value = "'11.99"
$ws.Range("E19").Value = "  -4.65%  "
$ws.Range("D20").Value = "'6.16"
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").Value = "'67.87"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "'235.43"
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("D27").Value = "'24.55"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("D30").Value = "'31.30"
$ws.Range("E30").Value = "  -5.79%  "
$ws.Range("D32").Value = "'139.77"
$ws.Range("E32").Value = "  -15.92%  "
$ws.Range("D34").Value = "'17.48"
$ws.Range("E34").Value = "  -3.08%  "
$ws.Range("D35").Value = "'0.0696"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("E37").Value = "  -4.46%  "
$ws.Range("E38").Value = "  +2.37%  "
$ws.Range("D39").Value = "'0.101"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "'22.39"
$ws.Range("E40").Value = "  +23.32%  "
$ws.Range("D41").Value = "'2.73"
$ws.Range("E41").Value = "  -2.57%  "
$ws.Range("E42").Value = "  -1.34%  "
$ws.Range("D43").Value = "'1.932.26"
$ws.Range("E43").Value = "  -3.14%  "
$ws.Range("D44").Value = "'0.0278"
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("D45").Value = "'10.23"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").Value = "'2.70"
$ws.Range("E47").Value = "  -2.77%  "
$ws.Range("D48").Value = "'2.88"
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("D49").Value = "'2.552.41"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").Value = "'52.65"
$ws.Range("E50").Value = "  -2.00%  "
$ws.Range("D51").Value = "'72.13"
